# Generate Report for Handoff
# b.md has now been handed off for localization: update status, handoff
# file names/dates, and (for zh-cn/de-de) the "content duplicate" /
# "error detail" columns to reflect the freshly generated handoff package.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-15 16:35:52"

# --- zh-cn sheet -------------------------------------------------------
# (leading "'" forces text, so "False" lands as a shared string like the
# rest of the Content Duplicate column, not as a native boolean)
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-15 16:35:47"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e389bb62f6334bf2afa44b0a57ca0d4a0330c03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b7988f4a2e363d7dee27078baf018d047b24c6a/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 39.14285714

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-08-15 16:35:52"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e389bb62f6334bf2afa44b0a57ca0d4a0330c03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b7988f4a2e363d7dee27078baf018d047b24c6a/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 39.14285714
